$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$textCells = @("D5","D6","D9","D10","D11","D12","D15","D19","D20","D21","D22","D23","D24","D28","D30","D31","D33","D34","D36","D37","D38","D39","D40","D41","D42","D43","D44","D45","D47","D48","D49","D50","D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '65.838.03'
$ws.Range("E2").Value = '  -1.15%  '

$ws.Range("D3").Value = '3.512.61'
$ws.Range("E3").Value = '  -1.47%  '

$ws.Range("E4").Value = '  -0.01%  '

$ws.Range("D5").Value = '596.25'
$ws.Range("E5").Value = '  -1.71%  '

$ws.Range("D6").Value = '143.26'
$ws.Range("E6").Value = '  -1.39%  '

$ws.Range("D7").Value = '3.511.38'
$ws.Range("E7").Value = '  -1.50%  '

$ws.Range("D9").Value = '0.500'
$ws.Range("E9").Value = '  +0.26%  '

$ws.Range("D10").Value = '0.134'
$ws.Range("E10").Value = '  -1.77%  '

$ws.Range("D11").Value = '7.66'
$ws.Range("E11").Value = '  -3.50%  '

$ws.Range("D12").Value = '0.404'
$ws.Range("E12").Value = '  -2.44%  '

$ws.Range("D13").Value = '4.103.74'
$ws.Range("E13").Value = '  -1.47%  '

$ws.Range("E14").Value = '  -3.64%  '

$ws.Range("D15").Value = '28.75'
$ws.Range("E15").Value = '  -4.49%  '

$ws.Range("D16").Value = '3.497.34'
$ws.Range("E16").Value = '  -1.86%  '

$ws.Range("E17").Value = '  +1.29%  '

$ws.Range("D18").Value = '65.762.77'
$ws.Range("E18").Value = '  -1.28%  '

$ws.Range("D19").Value = '10.90'
$ws.Range("E19").Value = '  -5.66%  '

$ws.Range("D20").Value = '6.19'
$ws.Range("E20").Value = '  -0.12%  '

$ws.Range("D21").Value = '14.39'
$ws.Range("E21").Value = '  -3.16%  '

$ws.Range("D22").Value = '413.17'
$ws.Range("E22").Value = '  -4.42%  '

$ws.Range("D23").Value = '0.596'
$ws.Range("E23").Value = '  -2.48%  '

$ws.Range("D24").Value = '77.49'
$ws.Range("E24").Value = '  -3.00%  '

$ws.Range("D25").Value = '3.648.12'
$ws.Range("E25").Value = '  -1.56%  '

$ws.Range("E26").Value = '  +0.08%  '

$ws.Range("E27").Value = '  -4.22%  '

$ws.Range("D28").Value = '9.03'
$ws.Range("E28").Value = '  -1.73%  '

$ws.Range("E29").Value = '  -3.28%  '

$ws.Range("D30").Value = '7.73'
$ws.Range("E30").Value = '  -3.67%  '

$ws.Range("D31").Value = '1.00'
$ws.Range("E31").Value = '  +0.55%  '

$ws.Range("D32").Value = '3.505.30'
$ws.Range("E32").Value = '  -1.45%  '

$ws.Range("D33").Value = '0.154'
$ws.Range("E33").Value = '  -0.22%  '

$ws.Range("D34").Value = '24.29'
$ws.Range("E34").Value = '  -4.31%  '

$ws.Range("E35").Value = '  +0.01%  '

$ws.Range("D36").Value = '7.48'
$ws.Range("E36").Value = '  -5.04%  '

$ws.Range("B37").Value = 'Fetch.AI'
$ws.Range("C37").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D37").Value = '1.26'
$ws.Range("E37").Value = '  -13.78%  '

$ws.Range("B38").Value = 'Monero'
$ws.Range("C38").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D38").Value = '175.48'
$ws.Range("E38").Value = '  +0.44%  '

$ws.Range("D39").Value = '5.21'
$ws.Range("E39").Value = '  -7.33%  '

$ws.Range("D40").Value = '1.58'
$ws.Range("E40").Value = '  -8.63%  '

$ws.Range("D41").Value = '0.0819'
$ws.Range("E41").Value = '  -3.72%  '

$ws.Range("D42").Value = '5.04'
$ws.Range("E42").Value = '  -3.09%  '

$ws.Range("D43").Value = '0.856'
$ws.Range("E43").Value = '  -3.79%  '

$ws.Range("D44").Value = '45.30'
$ws.Range("E44").Value = '  -1.84%  '

$ws.Range("D45").Value = '1.78'
$ws.Range("E45").Value = '  -8.38%  '

$ws.Range("E46").Value = '  +0.05%  '

$ws.Range("D47").Value = '2.41'
$ws.Range("E47").Value = '  -4.55%  '

$ws.Range("D48").Value = '7.08'
$ws.Range("E48").Value = '  -1.14%  '

$ws.Range("B49").Value = 'ONDO'
$ws.Range("C49").Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range("D49").Value = '1.08'
$ws.Range("E49").Value = '  -8.68%  '

$ws.Range("B50").Value = 'EnergySwap'
$ws.Range("C50").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D50").Value = '22.31'
$ws.Range("E50").Value = '  -5.95%  '

$ws.Range("D51").Value = '22.90'
$ws.Range("E51").Value = '  -8.84%  '
